$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2027.25
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2027.25
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6081.75
$ws.Range("N17").Value = -6417.75

$ws.Range("H18").Value = 910
$ws.Range("I18").Value = 910
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 910
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -626

$ws.Range("H33").Value = 3000471
$ws.Range("I33").Value = 4400208
$ws.Range("J33").Value = 1034.8572
$ws.Range("K33").Value = 4400208
$ws.Range("L33").Value = 1034.8572
$ws.Range("M33").Value = -4399979
$ws.Range("N33").Value = -1492.8572

$ws.Range("H41").Value = 856
$ws.Range("I41").Value = 1256.3
$ws.Range("J41").Value = 411.22223
$ws.Range("K41").Value = 1256.3
$ws.Range("L41").Value = 411.22223
$ws.Range("M41").Value = -816.3
$ws.Range("N41").Value = -1291.22223

$ws.Range("H53").Value = 1798.16
$ws.Range("I53").Value = 1548.25
$ws.Range("J53").Value = 2242.4443
$ws.Range("K53").Value = 1548.25
$ws.Range("L53").Value = 2242.4443
$ws.Range("M53").Value = -911.25
$ws.Range("N53").Value = -3516.4443

$ws.Range("H86").Value = 2344.4546
$ws.Range("I86").Value = 2365
$ws.Range("J86").Value = 2336.75
$ws.Range("K86").Value = 2365
$ws.Range("L86").Value = 2336.75
$ws.Range("M86").Value = -1242
$ws.Range("N86").Value = -4582.75

$ws.Range("H89").Value = 2344.4546
$ws.Range("I89").Value = 2365
$ws.Range("J89").Value = 2336.75
$ws.Range("K89").Value = 11825
$ws.Range("L89").Value = 11683.75
$ws.Range("M89").Value = -6209
$ws.Range("N89").Value = -22915.75

$ws.Range("H129").Value = 1564.7241
$ws.Range("I129").Value = 1231.909
$ws.Range("J129").Value = 1768.1111
$ws.Range("K129").Value = 3695.727
$ws.Range("L129").Value = 5304.3333
$ws.Range("M129").Value = 1304.273
$ws.Range("N129").Value = -15304.3333

$ws.Range("H137").Value = 3998.5
$ws.Range("I137").Value = 3998.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 11995.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -9445.5
$ws.Range("N137").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2215.5
$ws.Range("I2").Value = 995.4545000000001
$ws.Range("J2").Value = 3435.5454
$ws.Range("K2").Value = 995.4545000000001
$ws.Range("L2").Value = 3435.5454
$ws.Range("M2").Value = -882.4545000000001
$ws.Range("N2").Value = -3661.5454

$ws.Range("H32").Value = 1223.4916
$ws.Range("I32").Value = 1243.7413
$ws.Range("J32").Value = 49
$ws.Range("K32").Value = 1243.7413
$ws.Range("L32").Value = 49
$ws.Range("M32").Value = -956.7412999999999
$ws.Range("N32").Value = -623

$ws.Range("H61").Value = 3700.6553
$ws.Range("I61").Value = 3538.4614
$ws.Range("J61").Value = 5106.3335
$ws.Range("K61").Value = 3538.4614
$ws.Range("L61").Value = 5106.3335
$ws.Range("M61").Value = -3326.4614
$ws.Range("N61").Value = -5530.3335

$ws.Range("H116").Value = 2215.5
$ws.Range("I116").Value = 995.4545000000001
$ws.Range("J116").Value = 3435.5454
$ws.Range("K116").Value = 995.4545000000001
$ws.Range("L116").Value = 3435.5454
$ws.Range("M116").Value = 1298.5455
$ws.Range("N116").Value = -8023.5454

$ws.Range("H122").Value = 1111.625
$ws.Range("I122").Value = 897.8
$ws.Range("J122").Value = 1468
$ws.Range("K122").Value = 2693.4
$ws.Range("L122").Value = 4404
$ws.Range("M122").Value = -243.3999999999996
$ws.Range("N122").Value = -9304

$ws.Range("H136").Value = 3700.6553
$ws.Range("I136").Value = 3538.4614
$ws.Range("J136").Value = 5106.3335
$ws.Range("K136").Value = 10615.3842
$ws.Range("L136").Value = 15319.0005
$ws.Range("M136").Value = -8065.3842
$ws.Range("N136").Value = -20419.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2215.5
$ws.Range("I3").Value = 995.4545000000001
$ws.Range("J3").Value = 3435.5454
$ws.Range("K3").Value = 995.4545000000001
$ws.Range("L3").Value = 3435.5454
$ws.Range("M3").Value = -881.4545000000001
$ws.Range("N3").Value = -3663.5454

$ws.Range("H22").Value = 279.4
$ws.Range("I22").Value = 249.25
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 249.25
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = -76.25
$ws.Range("N22").Value = -746

$ws.Range("H94").Value = 788.875
$ws.Range("I94").Value = 746.5357
$ws.Range("J94").Value = 1085.25
$ws.Range("K94").Value = 746.5357
$ws.Range("L94").Value = 1085.25
$ws.Range("M94").Value = -295.5357
$ws.Range("N94").Value = -1987.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 999.6667
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 1049.5
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 1049.5
$ws.Range("M22").Value = -550
$ws.Range("N22").Value = -1749.5

$ws.Range("H31").Value = 3228.2856
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3228.2856
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3228.2856
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -3818.2856

$ws.Range("H34").Value = 3228.2856
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3228.2856
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3228.2856
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -3632.2856

$ws.Range("H58").Value = 1769.081
$ws.Range("I58").Value = 1526.75
$ws.Range("J58").Value = 3320
$ws.Range("K58").Value = 1526.75
$ws.Range("L58").Value = 3320
$ws.Range("M58").Value = -1323.75
$ws.Range("N58").Value = -3726

$ws.Range("H122").Value = 2097.75
$ws.Range("I122").Value = 1624.5
$ws.Range("J122").Value = 2571
$ws.Range("K122").Value = 4873.5
$ws.Range("L122").Value = 7713
$ws.Range("M122").Value = -2423.5
$ws.Range("N122").Value = -12613

$ws.Range("H132").Value = 3435.932
$ws.Range("I132").Value = 2815.2632
$ws.Range("J132").Value = 7366.8335
$ws.Range("K132").Value = 8445.7896
$ws.Range("L132").Value = 22100.5005
$ws.Range("M132").Value = -5915.7896
$ws.Range("N132").Value = -27160.5005

$ws.Range("H133").Value = 40000
$ws.Range("I133").Value = 40000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 40000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -37470
$ws.Range("N133").Value = ""

$ws.Range("H136").Value = 1769.081
$ws.Range("I136").Value = 1526.75
$ws.Range("J136").Value = 3320
$ws.Range("K136").Value = 4580.25
$ws.Range("L136").Value = 9960
$ws.Range("M136").Value = -2030.25
$ws.Range("N136").Value = -15060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4914.4165
$ws.Range("I134").Value = 5837.1665
$ws.Range("J134").Value = 3991.6667
$ws.Range("K134").Value = 17511.4995
$ws.Range("L134").Value = 11975.0001
$ws.Range("M134").Value = -12441.4995
$ws.Range("N134").Value = -22115.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 447499.5
$ws.Range("I20").Value = 887000
$ws.Range("J20").Value = 7999
$ws.Range("K20").Value = 887000
$ws.Range("L20").Value = 7999
$ws.Range("M20").Value = -886774
$ws.Range("N20").Value = -8451

$ws.Range("H22").Value = 1673.1875
$ws.Range("I22").Value = 1190.5385
$ws.Range("J22").Value = 3764.6667
$ws.Range("K22").Value = 1190.5385
$ws.Range("L22").Value = 3764.6667
$ws.Range("M22").Value = -895.5385000000001
$ws.Range("N22").Value = -4354.6667

$ws.Range("H27").Value = 1673.1875
$ws.Range("I27").Value = 1190.5385
$ws.Range("J27").Value = 3764.6667
$ws.Range("K27").Value = 1190.5385
$ws.Range("L27").Value = 3764.6667
$ws.Range("M27").Value = -1083.5385
$ws.Range("N27").Value = -3978.6667

$ws.Range("H40").Value = 9294.799999999999
$ws.Range("I40").Value = 7800
$ws.Range("J40").Value = 9668.5
$ws.Range("K40").Value = 7800
$ws.Range("L40").Value = 9668.5
$ws.Range("M40").Value = -7664
$ws.Range("N40").Value = -9940.5

$ws.Range("H55").Value = 508.55554
$ws.Range("I55").Value = 415.875
$ws.Range("J55").Value = 1250
$ws.Range("K55").Value = 415.875
$ws.Range("L55").Value = 1250
$ws.Range("M55").Value = -242.875
$ws.Range("N55").Value = -1596

$ws.Range("H68").Value = 1749.5
$ws.Range("I68").Value = 1499
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1499
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -750
$ws.Range("N68").Value = -3498

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = ""

$ws.Range("H71").Value = 1749.5
$ws.Range("I71").Value = 1499
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 7495
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -3751
$ws.Range("N71").Value = -17488

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = ""

$ws.Range("H136").Value = 27781260
$ws.Range("I136").Value = 3045.7666
$ws.Range("J136").Value = 166672340
$ws.Range("K136").Value = 9137.299800000001
$ws.Range("L136").Value = 500017020
$ws.Range("M136").Value = -6587.299800000001
$ws.Range("N136").Value = -500022120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 27000
$ws.Range("I62").Value = 50000
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 50000
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -49376
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 27000
$ws.Range("I65").Value = 50000
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 250000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -246880
$ws.Range("N65").Value = -26240

$ws.Range("H132").Value = 2181.5
$ws.Range("I132").Value = 2129.6667
$ws.Range("J132").Value = 2233.3333
$ws.Range("K132").Value = 6389.000100000001
$ws.Range("L132").Value = 6699.999899999999
$ws.Range("M132").Value = -3859.000100000001
$ws.Range("N132").Value = -11759.9999
